$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.633.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.777.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.82"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.79"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.775.45"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.16%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.99"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.410.89"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.768.15"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.49"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.606.73"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "458.83"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.24"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.99"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.08%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.22"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.56"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.55%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.991"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.76"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.25"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.03"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.297"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.82"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.28"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "393.07"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.77"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.33%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.71%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.721.95"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.81%  "
